# Auto-generated edit script: update crypto price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.918.24"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "2.099.00"
$ws.Range("E3").Value = "  +9.94%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'253.27"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'49.22"
$ws.Range("E8").Value = "  +5.53%  "
$ws.Range("D9").Value = "'60.75"
$ws.Range("E9").Value = "  +4.92%  "
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D12").Value = "'0.110"
$ws.Range("E12").Value = "  +10.24%  "
$ws.Range("D13").Value = "'14.87"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "2.403.73"
$ws.Range("E14").Value = "  +9.91%  "
$ws.Range("D15").Value = "'0.840"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("D16").Value = "2.115.56"
$ws.Range("E16").Value = "  +10.74%  "
$ws.Range("D17").Value = "'5.16"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "36.781.01"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'73.36"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").Value = "0.0₃0823"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").Value = "'13.37"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").Value = "'241.98"
$ws.Range("E22").Value = "  -3.40%  "
$ws.Range("D23").Value = "'5.34"
$ws.Range("E23").Value = "  +4.28%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'2.55"
$ws.Range("E25").Value = "  +2.47%  "
$ws.Range("D26").Value = "'171.82"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("D27").Value = "'9.47"
$ws.Range("E27").Value = "  +8.30%  "
$ws.Range("D28").Value = "'21.24"
$ws.Range("E28").Value = "  +14.01%  "
$ws.Range("E29").Value = "  -9.02%  "
$ws.Range("D30").Value = "'26.51"
$ws.Range("E30").Value = "  +40.91%  "
$ws.Range("E31").Value = "  -3.98%  "
$ws.Range("D32").Value = "'4.52"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").Value = "'0.0616"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "'1.02"
$ws.Range("E34").Value = "  +17.40%  "
$ws.Range("D35").Value = "'0.0925"
$ws.Range("E35").Value = "  +5.35%  "
$ws.Range("D36").Value = "'2.41"
$ws.Range("E36").Value = "  +21.36%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'1.86"
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("E39").Value = "  -4.41%  "
$ws.Range("D40").Value = "'1.34"
$ws.Range("E40").Value = "  -9.18%  "
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").Value = "'1.17"
$ws.Range("E42").Value = "  +7.82%  "
$ws.Range("D43").Value = "'98.37"
$ws.Range("D44").Value = "'16.79"
$ws.Range("E44").Value = "  -5.61%  "
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").Value = "1.343.97"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'0.0853"
$ws.Range("E47").Value = "  +4.62%  "
$ws.Range("D48").Value = "'7.14"
$ws.Range("E48").Value = "  +10.81%  "
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "2.280.74"
$ws.Range("E50").Value = "  +9.51%  "
$ws.Range("E51").Value = "  -3.99%  "
